$d = $word.ActiveDocument

# ===================================================================
# 1. First paragraph ("This is a Microsoft word document."):
#    - add two trailing spaces to the existing sentence
#    - append a new parenthetical note, in a dark-red color (C00000),
#      split across three runs exactly like the authored edit.
# ===================================================================
$p1 = $d.Paragraphs(1)
$p1Body = $p1.Range.Duplicate
$p1Body.MoveEnd(1, -1) | Out-Null          # exclude the paragraph mark
$p1Body.Text = "This is a Microsoft word document.  "

$insertPoint = $p1Body.Duplicate
$insertPoint.Collapse(0) | Out-Null        # wdCollapseEnd

$dash = [char]0x2013

$run1 = $d.Range($insertPoint.End, $insertPoint.End)
$run1.InsertAfter("(This is a change " + $dash + " Ve")
$run1.Font.Color = 192                     # C00000

$run2 = $d.Range($run1.End, $run1.End)
$run2.InsertAfter("rsion for branch alternate")
$run2.Font.Color = 192                     # C00000

$run3 = $d.Range($run2.End, $run2.End)
$run3.InsertAfter(")")
$run3.Font.Color = 192                     # C00000

# ===================================================================
# 2. The empty paragraph that used to sit just before "The Raven"
#    gets new paragraph-mark formatting: Calibri, bold, color 202122,
#    shaded F9F9F9. It stays empty (no visible text).
#
#    The runtime's Range.Font setters only take effect (and only sync
#    back onto the empty paragraph's mark / w:pPr/w:rPr) when the
#    range actually contains a character, so a throw-away placeholder
#    character is inserted, formatted, then removed again.
# ===================================================================
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text -eq [char]13 -and $candidate.Next().Range.Text -eq "The Raven" + [char]13) {
        $target = $candidate
        break
    }
}

if ($target -ne $null) {
    $target.Range.InsertBefore("X")
    $full = $target.Range

    $full.Font.Name = "Calibri"
    $full.Font.NameFarEast = "Times New Roman"
    $full.Font.NameBi = "Calibri"
    $full.Font.Bold = 1
    $full.Font.BoldBi = 1
    $full.Font.Color = 2236704             # 202122

    $full.Shading.Texture = 0              # wdTextureNone -> shd val="clear"
    $full.Shading.ForegroundPatternColor = -16777216  # wdColorAutomatic -> color="auto"
    $full.Shading.BackgroundPatternColor = 16382457   # F9F9F9

    $placeholder = $d.Range($target.Range.Start, $target.Range.Start + 1)
    $placeholder.Text = ""
}

# ===================================================================
# 3. Drop the trailing "ank God almighty, we are free at last."
#    paragraph's text, leaving an empty paragraph behind.
# ===================================================================
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastBody = $lastPara.Range.Duplicate
$lastBody.MoveEnd(1, -1) | Out-Null
$lastBody.Text = ""
